$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the account summary figures:
#  - VALOR MORA total
#  - Cant. Trabajadores (workers count) and Cant. Periodos (periods count)
$ws.Range("E11").Value = 2134400
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 40

# The detail table (rows 16-54) lists overdue periods for worker
# KAREN MARGARITA ABELLA SALGADO. Re-sequence the "Periodo Mora" column
# (E) from newest-first to oldest-first (ascending 2205 -> 2507), mirroring
# how the refreshed database now orders these records.
$periods = @(
    "2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
    "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
    "2501","2502","2503","2504","2505","2506","2507"
)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E" + $row).Value = $periods[$i]
}

# Row 55 previously held a single leftover record for a different worker
# (MARIA ALEJANDRA ESTRADA LOPEZ, period 2507). Replace it with the new
# "part 1" record: one more period (2508) for KAREN MARGARITA ABELLA SALGADO,
# continuing the same worker's history with her usual salary figures.
$ws.Range("C55").Value = "1047379783"
$ws.Range("D55").Value = "KAREN MARGARITA ABELLA SALGADO"
$ws.Range("E55").Value = "2508"
$ws.Range("F55").Value = 53360
$ws.Range("G55").Value = 1334000
